# The workbook gained a new weekly price record. It is inserted as a new
# row 13 on the "Espárragos" sheet, pushing every existing data row
# (13-56) down by one (to 14-57). The sheet's used range therefore grows
# from A1:R56 to A1:R57 (Excel maintains this automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 13, shifting rows 13:56 down
# to 14:57 (formatting carries down from the row below, matching Excel's
# default insert behaviour).
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with the new record.
$ws.Cells.Item(13, 1).Value  = 5
$ws.Cells.Item(13, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(13, 3).Value  = "Maule"
$ws.Cells.Item(13, 4).Value  = 44525
$ws.Cells.Item(13, 5).Value  = 7
$ws.Cells.Item(13, 6).Value  = 300000000
$ws.Cells.Item(13, 7).Value  = "Espárragos"
$ws.Cells.Item(13, 8).Value  = "Verde"
$ws.Cells.Item(13, 9).Value  = "Primera"
$ws.Cells.Item(13, 10).Value = 3000
$ws.Cells.Item(13, 11).Value = 1200
$ws.Cells.Item(13, 12).Value = 1200
$ws.Cells.Item(13, 13).Value = 1200
$ws.Cells.Item(13, 14).Value = "$/kilo"
$ws.Cells.Item(13, 15).Value = "Provincia de Linares"
$ws.Cells.Item(13, 16).Value = 1200
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(13, 18).Value = "Hortaliza"
